$p = $ppt.ActivePresentation

# 1. Table on slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") gets a new table style.
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{191D0C8E-A5B5-4785-8EC8-579F48F7C6C6}")
    }
}

# 2. The deck's theme (shared by every slide / the slide master) swaps its colour
#    scheme from the custom "Integral" palette back to the stock "Office Theme"
#    palette (the notes master keeps the palette that used to belong to the
#    slide master). Apply the swap through the presentation-wide theme colour
#    scheme so every slide (and therefore the underlying theme part) picks up
#    the new values.
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
